$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030144226804358
$ws.Range("D2").Value = 1.033091962115486
$ws.Range("E2").Value = 1.043816582810648
$ws.Range("F2").Value = 1.051627515709619
$ws.Range("I2").Value = 1.034050737202023
$ws.Range("J2").Value = 1.035287095779168
$ws.Range("K2").Value = 1.035895308166132
$ws.Range("L2").Value = 1.046589382286689
$ws.Range("M2").Value = 1.054378492767877
$ws.Range("N2").Value = 1.015753788384565
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031073561596111
$ws.Range("D3").Value = 1.033759860065317
$ws.Range("E3").Value = 1.044733880241246
$ws.Range("F3").Value = 1.052623478992939
$ws.Range("I3").Value = 1.034223132434526
$ws.Range("J3").Value = 1.035857863517911
$ws.Range("K3").Value = 1.036372491227186
$ws.Range("L3").Value = 1.047317514135319
$ws.Range("M3").Value = 1.055186668930673
$ws.Range("N3").Value = 1.015943837631545
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03167521008219
$ws.Range("D4").Value = 1.034191895735461
$ws.Range("E4").Value = 1.045328070694878
$ws.Range("F4").Value = 1.053268567374353
$ws.Range("I4").Value = 1.034332963048117
$ws.Range("J4").Value = 1.036226882647598
$ws.Range("K4").Value = 1.036680441559745
$ws.Range("L4").Value = 1.047788672304537
$ws.Range("M4").Value = 1.05570963197563
$ws.Range("N4").Value = 1.016066671619184
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031928215598251
$ws.Range("D5").Value = 1.034373488841803
$ws.Range("E5").Value = 1.045578019765598
$ws.Range("F5").Value = 1.05353991295588
$ws.Range("I5").Value = 1.034378723220846
$ws.Range("J5").Value = 1.036381944165567
$ws.Range("K5").Value = 1.036809706770654
$ws.Range("L5").Value = 1.047986748078594
$ws.Range("M5").Value = 1.055929489048679
$ws.Range("N5").Value = 1.016118277022869
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031970700561371
$ws.Range("D6").Value = 1.034403977047777
$ws.Range("E6").Value = 1.04561999615622
$ws.Range("F6").Value = 1.053585481857851
$ws.Range("I6").Value = 1.034386382338861
$ws.Range("J6").Value = 1.036407975303166
$ws.Range("K6").Value = 1.036831399379637
$ws.Range("L6").Value = 1.048020005892869
$ws.Range("M6").Value = 1.055966404185215
$ws.Range("N6").Value = 1.016126939793013
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031678590470607
$ws.Range("D7").Value = 1.034194322329537
$ws.Range("E7").Value = 1.045331409933261
$ws.Range("F7").Value = 1.053272192517787
$ws.Range("I7").Value = 1.034333576119966
$ws.Range("J7").Value = 1.036228954879611
$ws.Range("K7").Value = 1.036682169583631
$ws.Range("L7").Value = 1.047791319000623
$ws.Range("M7").Value = 1.055712569704648
$ws.Range("N7").Value = 1.016067361307168
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030458235772104
$ws.Range("D8").Value = 1.033317709275712
$ws.Range("E8").Value = 1.044126454981304
$ws.Range("F8").Value = 1.051963974238174
$ws.Range("I8").Value = 1.034109354828899
$ws.Range("J8").Value = 1.035480052164511
$ws.Range("K8").Value = 1.036056743377738
$ws.Range("L8").Value = 1.046835455514537
$ws.Range("M8").Value = 1.054651615058115
$ws.Range("N8").Value = 1.015818045316735
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028310205374158
$ws.Range("D9").Value = 1.031772003206075
$ws.Range("E9").Value = 1.04200810955577
$ws.Range("F9").Value = 1.049663633344227
$ws.Range("I9").Value = 1.033701100767489
$ws.Range("J9").Value = 1.034158088688227
$ws.Range("K9").Value = 1.034948435220992
$ws.Range("L9").Value = 1.045151211870699
$ws.Range("M9").Value = 1.052782273491563
$ws.Range("N9").Value = 1.015377657312128
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026879846371937
$ws.Range("D10").Value = 1.030740942127144
$ws.Range("E10").Value = 1.040599273343545
$ws.Range("F10").Value = 1.048133444088038
$ws.Range("I10").Value = 1.033420132157148
$ws.Range("J10").Value = 1.033275286197367
$ws.Range("K10").Value = 1.034205436530285
$ws.Range("L10").Value = 1.044028519340247
$ws.Range("M10").Value = 1.051536244159134
$ws.Range("N10").Value = 1.015083372042376
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026260890022489
$ws.Range("D11").Value = 1.030294358530885
$ws.Range("E11").Value = 1.039990051781232
$ws.Range("F11").Value = 1.047471669636734
$ws.Range("I11").Value = 1.033296390874525
$ws.Range("J11").Value = 1.032892679731684
$ws.Range("K11").Value = 1.033882743531065
$ws.Range("L11").Value = 1.043542426304587
$ws.Range("M11").Value = 1.050996760359195
$ws.Range("N11").Value = 1.014955783024868
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026031042762226
$ws.Range("D12").Value = 1.030128459724715
$ws.Range("E12").Value = 1.039763883202848
$ws.Range("F12").Value = 1.047225979608031
$ws.Range("I12").Value = 1.033250115977073
$ws.Range("J12").Value = 1.032750511342049
$ws.Range("K12").Value = 1.033762736259247
$ws.Range("L12").Value = 1.043361876741803
$ws.Range("M12").Value = 1.050796381381349
$ws.Range("N12").Value = 1.014908366835021
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02608034302054
$ws.Range("D13").Value = 1.030164046364386
$ws.Range("E13").Value = 1.039812391526195
$ws.Range("F13").Value = 1.047278675393038
$ws.Range("I13").Value = 1.033260056202694
$ws.Range("J13").Value = 1.032781009252701
$ws.Range("K13").Value = 1.033788484760169
$ws.Range("L13").Value = 1.043400604901756
$ws.Range("M13").Value = 1.050839362913531
$ws.Range("N13").Value = 1.014918538847203
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026241889533562
$ws.Range("D14").Value = 1.030280645638575
$ws.Range("E14").Value = 1.039971354091188
$ws.Range("F14").Value = 1.04745135831744
$ws.Range("I14").Value = 1.033292572135632
$ws.Range("J14").Value = 1.032880929090134
$ws.Range("K14").Value = 1.033872826639345
$ws.Range("L14").Value = 1.043527501872203
$ws.Range("M14").Value = 1.050980196774096
$ws.Range("N14").Value = 1.014951864074483
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026341431737844
$ws.Range("D15").Value = 1.0303524839849
$ws.Range("E15").Value = 1.040069312565243
$ws.Range("F15").Value = 1.047557770228589
$ws.Range("I15").Value = 1.033312564967167
$ws.Range("J15").Value = 1.032942486225539
$ws.Range("K15").Value = 1.033924773305963
$ws.Range("L15").Value = 1.043605688245432
$ws.Range("M15").Value = 1.051066970430688
$ws.Range("N15").Value = 1.01497239368323
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026920932932306
$ws.Range("D16").Value = 1.030770577842657
$ws.Range("E16").Value = 1.040639722618173
$ws.Range("F16").Value = 1.048177380990953
$ws.Range("I16").Value = 1.033428300688234
$ws.Range("J16").Value = 1.033300671286232
$ws.Range("K16").Value = 1.034226832253961
$ws.Range("L16").Value = 1.044060780666317
$ws.Range("M16").Value = 1.051572049172751
$ws.Range("N16").Value = 1.01509183632917
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027284545993556
$ws.Range("D17").Value = 1.031032803940245
$ws.Range("E17").Value = 1.04099774448147
$ws.Range("F17").Value = 1.048566263470673
$ws.Range("I17").Value = 1.033500342209296
$ws.Range("J17").Value = 1.03352525896342
$ws.Range("K17").Value = 1.034416046922204
$ws.Range("L17").Value = 1.044346259590069
$ws.Range("M17").Value = 1.05188888721957
$ws.Range("N17").Value = 1.015166716521684
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027496673738978
$ws.Range("D18").Value = 1.031185743608777
$ws.Range("E18").Value = 1.041206651020091
$ws.Range("F18").Value = 1.048793169886751
$ws.Range("I18").Value = 1.033542162076201
$ws.Range("J18").Value = 1.033656223603238
$ws.Range("K18").Value = 1.034526319014552
$ws.Range("L18").Value = 1.044512778432634
$ws.Range("M18").Value = 1.052073698740493
$ws.Range("N18").Value = 1.015210377261853
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027569010306703
$ws.Range("D19").Value = 1.031237889942565
$ws.Range("E19").Value = 1.041277896020833
$ws.Range("F19").Value = 1.04887055228013
$ws.Range("I19").Value = 1.033556387503187
$ws.Range("J19").Value = 1.033700873444171
$ws.Range("K19").Value = 1.034563903039205
$ws.Range("L19").Value = 1.044569557671918
$ws.Range("M19").Value = 1.052136715551642
$ws.Range("N19").Value = 1.015225261788492
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027245529775614
$ws.Range("D20").Value = 1.031004670832093
$ws.Range("E20").Value = 1.040959324001604
$ws.Range("F20").Value = 1.048524532000807
$ws.Range("I20").Value = 1.033492633594151
$ws.Range("J20").Value = 1.033501166295496
$ws.Range("K20").Value = 1.034395755655863
$ws.Range("L20").Value = 1.044315630001377
$ws.Range("M20").Value = 1.051854892945181
$ws.Range("N20").Value = 1.015158684197874
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026194316430141
$ws.Range("D21").Value = 1.030246310538272
$ws.Range("E21").Value = 1.039924540161033
$ws.Range("F21").Value = 1.047400504110377
$ws.Range("I21").Value = 1.033283005608113
$ws.Range("J21").Value = 1.032851506598157
$ws.Range("K21").Value = 1.033847994045519
$ws.Range("L21").Value = 1.043490133675129
$ws.Range("M21").Value = 1.05093872444653
$ws.Range("N21").Value = 1.01494205128414
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025533728395911
$ws.Range("D22").Value = 1.029769397153184
$ws.Range("E22").Value = 1.039274645227039
$ws.Range("F22").Value = 1.046694492287091
$ws.Range("I22").Value = 1.033149400064725
$ws.Range("J22").Value = 1.032442743298618
$ws.Range("K22").Value = 1.033502757723062
$ws.Range("L22").Value = 1.042971152555459
$ws.Range("M22").Value = 1.050362747421919
$ws.Range("N22").Value = 1.014805707177969
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025883885020487
$ws.Range("D23").Value = 1.03002222719472
$ws.Range("E23").Value = 1.039619098689241
$ws.Range("F23").Value = 1.047068694916951
$ws.Range("I23").Value = 1.033220397692141
$ws.Range("J23").Value = 1.032659464277563
$ws.Range("K23").Value = 1.033685853035575
$ws.Range("L23").Value = 1.043246270032986
$ws.Range("M23").Value = 1.050668078228797
$ws.Range("N23").Value = 1.014877998778747
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027263159406447
$ws.Range("D24").Value = 1.031017383011066
$ws.Range("E24").Value = 1.040976684320377
$ws.Range("F24").Value = 1.048543388415271
$ws.Range("I24").Value = 1.033496117405797
$ws.Range("J24").Value = 1.033512052838803
$ws.Range("K24").Value = 1.034404924695046
$ws.Range("L24").Value = 1.044329470182348
$ws.Range("M24").Value = 1.051870253478357
$ws.Range("N24").Value = 1.015162313707512
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028865233868087
$ws.Range("D25").Value = 1.032171715402228
$ws.Range("E25").Value = 1.042555159752256
$ws.Range("F25").Value = 1.050257737665538
$ws.Range("I25").Value = 1.03380819809979
$ws.Range("J25").Value = 1.034500115097309
$ws.Range("K25").Value = 1.035235691181832
$ws.Range("L25").Value = 1.045586609731162
$ws.Range("M25").Value = 1.053265512992975
$ws.Range("N25").Value = 1.015491631845898
